$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("D5").Value = "read probability"
$ws3.Range("E5").Value = "write probability"

$categories = @("dijkstra","FFT4096","FFT512","math","mp3enc","mpeg4","qsort","rj_enc","rj_dec")
$readVals   = @(0.85,0.54,0.54,0.56,0.5,0.61,0.53,0.81,0.81)

for ($i = 0; $i -lt $categories.Length; $i++) {
    $row = 6 + $i
    $ws3.Range("C" + $row).Value = $categories[$i]
    $ws3.Range("D" + $row).Value = $readVals[$i]
}

$ws3.Range("E6").Formula = "=1-D6"
$ws3.Range("E7:E14").Formula = "=1-D7"

Write-Host "done"
